$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("RV32Z")

# Un-merge the cells that get split apart by this edit (header cells that
# used to span two rows, and the banner row at the bottom).
$ws.Range("A2:A3").UnMerge()
$ws.Range("B2:B3").UnMerge()
$ws.Range("A5:I5").UnMerge()

# Rows 3 and 4 become completely blank (no left-over formatting) once the
# merged header area is split apart.
$ws.Range("A3:I4").Clear()

# Fill in the coverage link for the FENCE.I requirement row.
$ws.Range("I2").Value = "isacov_agent.cov_model.rv32zifencei_fence_i_cg"
$ws.Range("I2").HorizontalAlignment = -4131
$ws.Range("I2").VerticalAlignment = -4160
$ws.Range("I2").WrapText = $true

# Add the new "Missing Coverage" column with its header and explanatory note.
$ws.Range("J1").Value = "Missing Coverage"
$ws.Range("J1").Font.Name = "DejaVu Sans"
$ws.Range("J1").Font.Bold = $true
$ws.Range("J1").Font.Size = 11

$ws.Range("J2").Value = "This only tracks that the instruction was executed. Refer to the vplans for the core under test for specific instruction test and coverage."
$ws.Range("J2").Font.Name = "DejaVu Sans"
$ws.Range("J2").Font.Bold = $false
$ws.Range("J2").Font.Size = 11
$ws.Range("J2").VerticalAlignment = -4160
$ws.Range("J2").WrapText = $true

$ws.Columns.Item(10).ColumnWidth = 31.6

$ws.Rows.Item(1).RowHeight = 27.6
$ws.Rows.Item(2).RowHeight = 89.4

$wb.Save()
